# "Add references and more"
#
# The "Variance rate" row (original row 3) is removed from the table; the
# "Modified rate" row (original row 4) moves up to row 3 and its label is
# expanded to "Modified variance rate". The previously-blank last row
# (original row 5) becomes row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Variance rate" row (row 3). This shifts the "Modified rate"
# row up from row 4 to row 3, and the trailing blank row up from row 5 to
# row 4 - exactly matching the target layout.
$ws.Rows("3").Delete() | Out-Null

# Update the (now row 3) label to reflect the merged/renamed field.
$ws.Range("A3").Value = 'Modified variance rate $( \tau_l )$'

# Reflect the row selection left after the edit (row 3 selected).
$ws.Rows("3").Select() | Out-Null
